# Fix typo in cell B7: "préférance" -> "préférence"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("B7").Value = "Ordre de préférence des stages de l'étudiant"

# Match the cursor/selection position reflected in the saved file
$ws.Range("C4").Select()
